# "replace ~ with 0"
# The diagonal cells of the distance matrix (distance from a place to
# itself) were stored as the text placeholder "~". Replace each of them
# with the number 0. Excel will then drop the now-unused "~" entry from
# the shared strings table automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$diagonalCells = @("B2", "C3", "D4", "E5", "F6", "G7", "H8", "I9", "J10", "K11")

foreach ($addr in $diagonalCells) {
    $ws.Range($addr).Value = 0
}

# Leave the selection on the last edited cell, matching the author's
# final state after making this change.
$ws.Range("K11").Select() | Out-Null
